$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.58102533333333
$ws.Range("H2").Value = 58.743076
$ws.Range("I2").Value = 0.3529199051285138
$ws.Range("J2").Value = 0.3529199051285138
$ws.Range("M2").Value = 46.63275166666667
$ws.Range("N2").Value = 139.898255
$ws.Range("O2").Value = 0.9158911059585902
$ws.Range("P2").Value = 0.9158911059585902
$ws.Range("Q2").Value = 913.1170917480423
$ws.Range("R2").Value = 8218.053825732381
$ws.Range("S2").Value = 0.3232362022229552
$ws.Range("T2").Value = 0.3232362022229553

$ws.Range("G3").Value = 19.58102533333333
$ws.Range("H3").Value = 58.743076
$ws.Range("I3").Value = 0.3529199051285138
$ws.Range("J3").Value = 0.3529199051285138
$ws.Range("M3").Value = 2.770761666666667
$ws.Range("N3").Value = 8.312284999999999
$ws.Range("O3").Value = 0.05441917700612491
$ws.Range("P3").Value = 0.05441917700612491
$ws.Range("Q3").Value = 54.25435438762889
$ws.Range("R3").Value = 488.28918948866
$ws.Range("S3").Value = 0.0192056107861734
$ws.Range("T3").Value = 0.01920561078617341

$ws.Range("G4").Value = 19.58102533333333
$ws.Range("H4").Value = 58.743076
$ws.Range("I4").Value = 0.3529199051285138
$ws.Range("J4").Value = 0.3529199051285138
$ws.Range("M4").Value = 0.8496050000000001
$ws.Range("N4").Value = 2.548815
$ws.Range("O4").Value = 0.01668667696558362
$ws.Range("P4").Value = 0.01668667696558362
$ws.Range("Q4").Value = 16.63613702832667
$ws.Range("R4").Value = 149.72523325494
$ws.Range("S4").Value = 0.005889060451603929
$ws.Range("T4").Value = 0.00588906045160393

$ws.Range("G5").Value = 19.58102533333333
$ws.Range("H5").Value = 58.743076
$ws.Range("I5").Value = 0.3529199051285138
$ws.Range("J5").Value = 0.3529199051285138
$ws.Range("M5").Value = 0.662052
$ws.Range("N5").Value = 1.986156
$ws.Range("O5").Value = 0.01300304006970129
$ws.Range("P5").Value = 0.0130030400697013
$ws.Range("Q5").Value = 12.963656983984
$ws.Range("R5").Value = 116.672912855856
$ws.Range("S5").Value = 0.004589031667781244
$ws.Range("T5").Value = 0.004589031667781246

$ws.Range("G6").Value = 14.68975
$ws.Range("H6").Value = 44.06925
$ws.Range("I6").Value = 0.2647616806631773
$ws.Range("J6").Value = 0.2647616806631773
$ws.Range("M6").Value = 46.63275166666667
$ws.Range("N6").Value = 139.898255
$ws.Range("O6").Value = 0.9158911059585902
$ws.Range("P6").Value = 0.9158911059585902
$ws.Range("Q6").Value = 685.0234637954167
$ws.Range("R6").Value = 6165.21117415875
$ws.Range("S6").Value = 0.2424928685180525
$ws.Range("T6").Value = 0.2424928685180526

$ws.Range("G7").Value = 14.68975
$ws.Range("H7").Value = 44.06925
$ws.Range("I7").Value = 0.2647616806631773
$ws.Range("J7").Value = 0.2647616806631773
$ws.Range("M7").Value = 2.770761666666667
$ws.Range("N7").Value = 8.312284999999999
$ws.Range("O7").Value = 0.05441917700612491
$ws.Range("P7").Value = 0.05441917700612491
$ws.Range("Q7").Value = 40.70179619291666
$ws.Range("R7").Value = 366.31616573625
$ws.Range("S7").Value = 0.01440811276444856
$ws.Range("T7").Value = 0.01440811276444857

$ws.Range("G8").Value = 14.68975
$ws.Range("H8").Value = 44.06925
$ws.Range("I8").Value = 0.2647616806631773
$ws.Range("J8").Value = 0.2647616806631773
$ws.Range("M8").Value = 0.8496050000000001
$ws.Range("N8").Value = 2.548815
$ws.Range("O8").Value = 0.01668667696558362
$ws.Range("P8").Value = 0.01668667696558362
$ws.Range("Q8").Value = 12.48048504875
$ws.Range("R8").Value = 112.32436543875
$ws.Range("S8").Value = 0.004417992638091448
$ws.Range("T8").Value = 0.004417992638091449

$ws.Range("G9").Value = 14.68975
$ws.Range("H9").Value = 44.06925
$ws.Range("I9").Value = 0.2647616806631773
$ws.Range("J9").Value = 0.2647616806631773
$ws.Range("M9").Value = 0.662052
$ws.Range("N9").Value = 1.986156
$ws.Range("O9").Value = 0.01300304006970129
$ws.Range("P9").Value = 0.0130030400697013
$ws.Range("Q9").Value = 9.725378366999999
$ws.Range("R9").Value = 87.528405303
$ws.Range("S9").Value = 0.003442706742584753
$ws.Range("T9").Value = 0.003442706742584754

$ws.Range("G10").Value = 19.14352733333333
$ws.Range("H10").Value = 57.430582
$ws.Range("I10").Value = 0.3450346309906436
$ws.Range("J10").Value = 0.3450346309906436
$ws.Range("M10").Value = 46.63275166666667
$ws.Range("N10").Value = 139.898255
$ws.Range("O10").Value = 0.9158911059585902
$ws.Range("P10").Value = 0.9158911059585902
$ws.Range("Q10").Value = 892.7153561593791
$ws.Range("R10").Value = 8034.438205434411
$ws.Range("S10").Value = 0.3160141497720347
$ws.Range("T10").Value = 0.3160141497720347

$ws.Range("G11").Value = 19.14352733333333
$ws.Range("H11").Value = 57.430582
$ws.Range("I11").Value = 0.3450346309906436
$ws.Range("J11").Value = 0.3450346309906436
$ws.Range("M11").Value = 2.770761666666667
$ws.Range("N11").Value = 8.312284999999999
$ws.Range("O11").Value = 0.05441917700612491
$ws.Range("P11").Value = 0.05441917700612491
$ws.Range("Q11").Value = 53.04215169998556
$ws.Range("R11").Value = 477.37936529987
$ws.Range("S11").Value = 0.01877650065712283
$ws.Range("T11").Value = 0.01877650065712283

$ws.Range("G12").Value = 19.14352733333333
$ws.Range("H12").Value = 57.430582
$ws.Range("I12").Value = 0.3450346309906436
$ws.Range("J12").Value = 0.3450346309906436
$ws.Range("M12").Value = 0.8496050000000001
$ws.Range("N12").Value = 2.548815
$ws.Range("O12").Value = 0.01668667696558362
$ws.Range("P12").Value = 0.01668667696558362
$ws.Range("Q12").Value = 16.26443654003667
$ws.Range("R12").Value = 146.37992886033
$ws.Range("S12").Value = 0.005757481429280219
$ws.Range("T12").Value = 0.005757481429280219

$ws.Range("G13").Value = 19.14352733333333
$ws.Range("H13").Value = 57.430582
$ws.Range("I13").Value = 0.3450346309906436
$ws.Range("J13").Value = 0.3450346309906436
$ws.Range("M13").Value = 0.662052
$ws.Range("N13").Value = 1.986156
$ws.Range("O13").Value = 0.01300304006970129
$ws.Range("P13").Value = 0.0130030400697013
$ws.Range("Q13").Value = 12.674010558088
$ws.Range("R13").Value = 114.066095022792
$ws.Range("S13").Value = 0.004486499132205939
$ws.Range("T13").Value = 0.004486499132205939

$ws.Range("G14").Value = 2.068613
$ws.Range("H14").Value = 6.205839
$ws.Range("I14").Value = 0.03728378321766519
$ws.Range("J14").Value = 0.0372837832176652
$ws.Range("M14").Value = 46.63275166666667
$ws.Range("N14").Value = 139.898255
$ws.Range("O14").Value = 0.9158911059585902
$ws.Range("P14").Value = 0.9158911059585902
$ws.Range("Q14").Value = 96.46511632343834
$ws.Range("R14").Value = 868.186046910945
$ws.Range("S14").Value = 0.0341478854455477
$ws.Range("T14").Value = 0.03414788544554771

$ws.Range("G15").Value = 2.068613
$ws.Range("H15").Value = 6.205839
$ws.Range("I15").Value = 0.03728378321766519
$ws.Range("J15").Value = 0.0372837832176652
$ws.Range("M15").Value = 2.770761666666667
$ws.Range("N15").Value = 8.312284999999999
$ws.Range("O15").Value = 0.05441917700612491
$ws.Range("P15").Value = 0.05441917700612491
$ws.Range("Q15").Value = 5.731633603568334
$ws.Range("R15").Value = 51.58470243211499
$ws.Range("S15").Value = 0.002028952798380111
$ws.Range("T15").Value = 0.002028952798380112

$ws.Range("G16").Value = 2.068613
$ws.Range("H16").Value = 6.205839
$ws.Range("I16").Value = 0.03728378321766519
$ws.Range("J16").Value = 0.0372837832176652
$ws.Range("M16").Value = 0.8496050000000001
$ws.Range("N16").Value = 2.548815
$ws.Range("O16").Value = 0.01668667696558362
$ws.Range("P16").Value = 0.01668667696558362
$ws.Range("Q16").Value = 1.757503947865
$ws.Range("R16").Value = 15.817535530785
$ws.Range("S16").Value = 0.000622142446608027
$ws.Range("T16").Value = 0.0006221424466080272

$ws.Range("G17").Value = 2.068613
$ws.Range("H17").Value = 6.205839
$ws.Range("I17").Value = 0.03728378321766519
$ws.Range("J17").Value = 0.0372837832176652
$ws.Range("M17").Value = 0.662052
$ws.Range("N17").Value = 1.986156
$ws.Range("O17").Value = 0.01300304006970129
$ws.Range("P17").Value = 0.0130030400697013
$ws.Range("Q17").Value = 1.369529373876
$ws.Range("R17").Value = 12.325764364884
$ws.Range("S17").Value = 0.0004848025271293571
$ws.Range("T17").Value = 0.0004848025271293573
